$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# A new weekly price record was added to the data set. In the source data
# it lands right after the existing row for this market/product (row 16),
# so insert a fresh row there: this pushes the previous rows 16-40 down to
# 17-41 (and grows the used range from A1:R40 to A1:R41).
$ws.Rows("16").Insert()

# Seed the newly inserted row 16 with a duplicate of the row right below it
# (the old row 16, now shifted to row 17) so every column - Mercado ID,
# Mercado, Region, Codreg, Categoria, Variedad, Calidad, precios, unidad,
# Origen, etc. - carries over with matching formatting/types.
$ws.Range("A16:R16").Value2 = $ws.Range("A17:R17").Value2

# Then apply the two values that actually differ for this new record:
# the reporting date and the traded volume.
$ws.Range("D16").Value2 = 44883
$ws.Range("J16").Value2 = 9700
